$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Sort the existing 4 product rows (A2:J5) ascending by the
#        "consecutive" column (C), like the author did before adding the
#        new rows. This reproduces the row reorder, the <sortState> element,
#        and keeps the already-clean numeric literals in column E untouched. ---
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("C2:C5"))
$sort.SetRange($ws.Range("A2:J5"))
$sort.Header = -4142
$sort.Apply()

# Column A has no data in this sheet; the Sort engine stamps empty styled
# cells across the whole sorted range, so drop the stray A2:A5 cells again.
$ws.Range("A2:A5").Clear()

# --- 2. Renumber the product codes (114001-114004 -> 14001-14004) now that
#        the rows are sorted into their final order. ---
$ws.Range("B2").Value = 14001
$ws.Range("B3").Value = 14002
$ws.Range("B4").Value = 14003
$ws.Range("B5").Value = 14004

# Row 5 ("Fadeli lija pligo agua #80") used the lowercase "unidad" unit;
# the refreshed sheet uses the same "Unidad" spelling as the other rows.
$ws.Range("F5").Value = "Unidad"

# --- 3. Append the three new products (rows 6-8). ---
$ws.Range("B6").Value = 14005
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = "Pintura spray negro mate"
$ws.Range("E6").Value = 2367.26
$ws.Range("F6").Value = "Unidad"
$ws.Range("G6").Value = "t"
$ws.Range("H6").Value = 13
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 14

$ws.Range("B7").Value = 14006
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = "Thinner corriente litro"
$ws.Range("E7").Value = 1876.11
$ws.Range("F7").Value = "Unidad"
$ws.Range("G7").Value = "t"
$ws.Range("H7").Value = 13
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 14

$ws.Range("B8").Value = 14007
$ws.Range("C8").Value = 7
$ws.Range("D8").Value = "Corrostop 9000-700-14 negro cuarto"
$ws.Range("E8").Value = 6402.65
$ws.Range("F8").Value = "Unidad"
$ws.Range("G8").Value = "t"
$ws.Range("H8").Value = 13
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 14

# --- 4. New "predetermined discount" column (L) for the new products. ---
$ws.Range("L6").Value = 5
$ws.Range("L6").Font.Size = 12
$ws.Range("L7").Value = 5
$ws.Range("L8").Value = 10
$ws.Range("L8").Font.Size = 12

# --- 5. Misc view/print bookkeeping the author's session picked up. ---
$ws.Range("C9").Select()
$ws.PageSetup.Orientation = 1
